# Apply cryptos list price/volume updates for Thu Oct 24 19:37:24 UTC 2024 run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.093.31"
$ws.Range("E2").Value = "  +2.81%  "
$ws.Range("D3").Value = "2.535.81"
$ws.Range("E3").Value = "  +1.45%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +2.53%  "
$ws.Range("D9").Value = "2.532.78"
$ws.Range("E9").Value = "  +1.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.142"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.164"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.16"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.24%  "
$ws.Range("D15").Value = "2.994.95"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "67.954.61"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").Value = "2.530.90"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("E19").Value = "  +5.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "364.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.04"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.59%  "
$ws.Range("D28").Value = "2.664.53"
$ws.Range("E28").Value = "  +1.53%  "
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").Value = "0.0₃0998"
$ws.Range("E30").Value = "  +3.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "549.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.79%  "
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.15%  "
$ws.Range("E35").Value = "  -0.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.60"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("E40").Value = "  +1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.358"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.58%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.09%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.565"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "147.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.29%  "
$ws.Range("D49").Value = "0.0₆0278"
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0758"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.00%  "
